# Append a new data row (row 17) to the sheet, mirroring the existing rows'
# layout: date/time in column A (inherits the existing date-time style from
# the column), numeric metrics in B:M, and the shared "Noun" method label in N.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 17

$ws.Range("A$row").Value = 42625.884247685186
$ws.Range("B$row").Value = -22
$ws.Range("C$row").Value = 60
$ws.Range("D$row").Value = 38
$ws.Range("E$row").Value = 61
$ws.Range("F$row").Value = 38
$ws.Range("G$row").Value = 8540
$ws.Range("H$row").Value = 7203
$ws.Range("I$row").Value = 1026
$ws.Range("J$row").Value = 174
$ws.Range("K$row").Value = 110
$ws.Range("L$row").Value = 8
$ws.Range("M$row").Value = 5
$ws.Range("N$row").Value = "Noun"
